$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 101; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $txt = $cell.Text
    if ($txt -ne $null -and $txt.Length -gt 0) {
        $newtxt = $txt.Replace("<b>", "").Replace("</b>", "")
        if ($newtxt -ne $txt) {
            $cell.Value = $newtxt
        }
    }
}
